# Procedures.xlsx update — "Actualizo excel de procedures y creo el de triggers"
#
# Substance of the edit (per the OOXML diff): five existing notes/signatures
# get their text extended/clarified, and one previously-empty note cell gets
# new text. Everything else in the sheet is unchanged (the rest of the diff
# noise is just shared-string re-indexing caused by these edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 — login(username, password) / Info note
$ws.Range("C8").Value = "Recordar que luego de 3 intentos fallidos se debe deshabilitar al usuario. Si el login es satisfactorio, limpiar intentos_login."

# Row 11 — modificarAfiliado(...) gains a new "motivo" parameter, and its
# note explains it is optional unless the plan changes
$ws.Range("A11").Value = "modificarAfiliado(password, direccion, telefono, mail, sexo, estadoCivil, familiaresACargo, idPlanMedico, motivo)"
$ws.Range("C11").Value = "Recordar que si se cambia el plan hay que agregar al anterior al historial de cambios de plan. Motivo es opcional si se cambia el plan."

# Row 13 — eliminarAfiliado(idAfiliado) / Info note
$ws.Range("C13").Value = "Baja lógica (inhabilitar). Recordar dar de baja los turnos que tenía posteriores a la fecha de inhabilitación para que otros los puedan usar."

# Row 18 — comprarBonos(...) previously had no note; now it does
$ws.Range("C18").Value = "Validar que el afiliado que los compra esté habilitado."

# Row 25 — registrarLlegada(idAfiliado, numeroTurno, fecha) / Info note
$ws.Range("C25").Value = "Crear una consulta médica para ese turno."

# Update the view state to match the author's final scroll/selection position
$ws.Range("C26").Select()
